$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.148.73'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.672.75'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.80%  '
$ws.Range('D5').Value = "'211.00"
$ws.Range('E5').Value = '  -3.68%  '
$ws.Range('D6').Value = "'0.5259"
$ws.Range('E6').Value = '  -4.22%  '
$ws.Range('E7').Value = '  -0.79%  '
$ws.Range('D8').Value = "'0.2655"
$ws.Range('E8').Value = '  -3.31%  '
$ws.Range('D9').Value = "'0.06297"
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('D10').Value = "'21.24"
$ws.Range('E10').Value = '  -3.53%  '
$ws.Range('D11').Value = "'0.07531"
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.682.51'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.448"
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').Value = "'0.5641"
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('D15').Value = "'0.000008018"
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').Value = "'66.49"
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '26.212.19'
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = "'4.811"
$ws.Range('E19').Value = '  -2.51%  '
$ws.Range('D20').Value = "'188.00"
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = "'10.40"
$ws.Range('E21').Value = '  -5.30%  '
$ws.Range('D22').Value = "'6.180"
$ws.Range('E22').Value = '  -1.13%  '
$ws.Range('D23').Value = "'1.003"
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('D24').Value = "'148.35"
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('D25').Value = "'0.1251"
$ws.Range('E25').Value = '  -5.76%  '
$ws.Range('D26').Value = "'7.592"
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('D28').Value = "'0.06217"
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').Value = "'1.356"
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').Value = "'1.279"
$ws.Range('E30').Value = '  -4.07%  '
$ws.Range('D31').Value = "'3.487"
$ws.Range('E31').Value = '  -3.49%  '
$ws.Range('D32').Value = "'3.437"
$ws.Range('E32').Value = '  -4.66%  '
$ws.Range('D33').Value = "'1.630"
$ws.Range('E33').Value = '  -3.19%  '
$ws.Range('D34').Value = "'1.001"
$ws.Range('E34').Value = '  -4.03%  '
$ws.Range('D35').Value = "'0.6052"
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('D36').Value = "'2.402"
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').Value = "'2.714"
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').Value = "'0.01615"
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('D40').Value = '1.076.92'
$ws.Range('E40').Value = '  -3.74%  '
$ws.Range('D41').Value = "'0.8663"
$ws.Range('E41').Value = '  -2.08%  '
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').Value = "'100.09"
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('D45').Value = "'0.00000000108"
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').Value = "'56.18"
$ws.Range('E46').Value = '  -2.33%  '
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').Value = "'7.993"
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('D49').Value = "'0.05239"
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('D50').Value = "'0.4254"
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('D51').Value = "'5.984"
$ws.Range('E51').Value = '  -2.15%  '
